# Arbeitszeit Topeiner - add new work-log entries (stream UI work)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data: E = date serial, F = hours, G = description
$newRows = @(
    @{ Row = 33; Date = 43716; Hours = 0.2;  Text = "Funktionalität der UI programmiert" },
    @{ Row = 34; Date = 43717; Hours = 0.3;  Text = "Kurzes Team-Meeting" },
    @{ Row = 35; Date = 43718; Hours = 0.5;  Text = "Statusupdate erstellt" },
    @{ Row = 36; Date = 43719; Hours = 0.75; Text = "Team-Meeting" },
    @{ Row = 37; Date = 43721; Hours = 4;    Text = "Pflichtenheft erstellt, Termine und andere Informationen für die Diplomarbeit erhalten" },
    @{ Row = 38; Date = 43724; Hours = 0.5;  Text = "Projektplan mit MS Project erstellt" },
    @{ Row = 39; Date = 43726; Hours = 1;    Text = "Projektplan mit Informationen meines Kollegen ergänzt und fertiggestellt" },
    @{ Row = 40; Date = 43732; Hours = 0.75; Text = "Statusupdate erstellt" },
    @{ Row = 41; Date = 43734; Hours = 0.5;  Text = "Pflichtenheft ergänzt" },
    @{ Row = 42; Date = 43735; Hours = 4;    Text = "Pflichtenheft fertiggestellt, Arbeiten an UI" },
    @{ Row = 43; Date = 43738; Hours = 2;    Text = "Diplomarbeitspräsentation erstellt" },
    @{ Row = 44; Date = 43758; Hours = 4;    Text = "Anpassungen an der UI vorgenommen (Aufbau der Funktion `"Hinzufügen eines Torrents`" verändert), Logo implementiert" },
    @{ Row = 45; Date = 43759; Hours = 1.5;  Text = "Sprache der UI geändert, Icon erstellt und eingebunden" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 5).Value = $entry.Date
    $ws.Cells.Item($r, 6).Value = $entry.Hours
    $ws.Cells.Item($r, 7).Value = $entry.Text
}

# Copy formatting (date format, alignment, etc.) down from the last existing data row
$fmtSource = $ws.Range("E32:G32")
$fmtTarget = $ws.Range("E33:G45")
$fmtSource.Copy()
$fmtTarget.PasteSpecial(-4122)

# Move selection / view to the new bottom of the list, as in the saved workbook
$ws.Cells.Item(46, 7).Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
